# Change title color of drawings
#
# - Slide 1 "TextBox 26"  (Smart contract PepitoDisguise ... title)
# - Slide 2 "TextBox 3"   (Pepito - React calling sequence title)
#   both get: fixed/enlarged box, light-green fill (E1EFD8), bold
#   dark-red (980000) text, centered vertically, no autofit.
# - Slide 2 gains the round logo picture that already exists on slide 1
#   ("Google Shape;117;p28"), copied across so it keeps the same embedded
#   image relationship.

$p  = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

$fillColor = 14217185   # RGB(0xE1,0xEF,0xD8) -> BGR-int used by VBA's RGB()
$fontColor = 152         # RGB(0x98,0x00,0x00) -> BGR-int used by VBA's RGB()

# ---------------------------------------------------------------------
# Slide 1: "TextBox 26" -> "Smart contract PepitoDisguise is inspired
# from MetaCoin but not exactly the same"
# ---------------------------------------------------------------------
$sh1 = $s1.Shapes.Item("TextBox 26")

$sh1.Left   = 36.61677169799805
$sh1.Top    = 21.285669326782227
$sh1.Width  = 682.4702758789062
$sh1.Height = 39.82133865356445

$sh1.Fill.ForeColor.RGB = $fillColor

$sh1.TextFrame.AutoSize       = 0   # ppAutoSizeNone
$sh1.TextFrame.VerticalAnchor = 3   # msoAnchorMiddle

$tr1 = $sh1.TextFrame.TextRange
# fix the "Samrt" -> "Smart" typo without disturbing the other two runs
$tr1.Characters(1, 15).Text = "Smart contract "
$tr1.Font.Bold      = $true
$tr1.Font.Color.RGB = $fontColor

# ---------------------------------------------------------------------
# Slide 2: "TextBox 3" -> "Pepito – React calling sequence"
# ---------------------------------------------------------------------
$sh2 = $s2.Shapes.Item("TextBox 3")

$sh2.Left   = 48.8571662902832
$sh2.Top    = 32.57141876220703
$sh2.Width  = 264.7867126464844
$sh2.Height = 46.79133987426758

$sh2.Fill.ForeColor.RGB = $fillColor

$sh2.TextFrame.AutoSize       = 0   # ppAutoSizeNone
$sh2.TextFrame.VerticalAnchor = 3   # msoAnchorMiddle

$tr2 = $sh2.TextFrame.TextRange
$tr2.Font.Bold      = $true
$tr2.Font.Color.RGB = $fontColor

# ---------------------------------------------------------------------
# Slide 2: add the round logo picture (copy of the one already on slide 1
# so the embedded image part is reused rather than duplicated).
# ---------------------------------------------------------------------
$logo = $s1.Shapes.Item("Google Shape;117;p28")
$logo.Copy()
$s2.Shapes.Paste() | Out-Null
